$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.306.85"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.854.31"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.41"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4557"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3894"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.04%  "
$ws.Range("E9").Value = "  -9.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07923"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.40%  "
$ws.Range("E11").Value = "  -2.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.38"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.864.31"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.913"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.165"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06631"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "85.99"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001026"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.21"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.08%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.504"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.317.78"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.88"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.89%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.291"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.083.94"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.22"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.93"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.064"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.473"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.25"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9475"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09353"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.447"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.588"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.257"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06034"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02228"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("E39").Value = "  -0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.065"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.45%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5930"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -2.98%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1885"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.18"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.282"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5608"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.10"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.386"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.23%  "
$ws.Range("E49").Value = "  -5.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06734"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.07"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.13%  "
